$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 29 - this shifts the existing rows 29-42 down to 30-43
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly record
$ws.Cells.Item(29, 1).Value = 7
$ws.Cells.Item(29, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(29, 3).Value = "Ñuble"
$ws.Cells.Item(29, 4).Value = 45027
$ws.Cells.Item(29, 5).Value = 16
$ws.Cells.Item(29, 6).Value = 100112044
$ws.Cells.Item(29, 7).Value = "Perejil"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 200
$ws.Cells.Item(29, 11).Value = 1500
$ws.Cells.Item(29, 12).Value = 1500
$ws.Cells.Item(29, 13).Value = 1500
$ws.Cells.Item(29, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(29, 15).Value = "Región del Maule"
$ws.Cells.Item(29, 16).Value = 1500
$ws.Cells.Item(29, 17).Value = 1
$ws.Cells.Item(29, 18).Value = "Hortaliza"
